$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, using the same style as the existing headers (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill columns I and J for data rows 2..38:
#   I = 1 (constant)
#   J = value copied from column H
for ($r = 2; $r -le 38; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
